# ---------------------------------------------------------------------------
# RF-RNF.xlsx - "Unificación de Metricas con RF-RNF"
#
# Adds three new worksheets (Metricas, MTBF, MTTR) at the end of the
# workbook containing defect-density / MTBF / MTTR metric tables, and fixes
# the priority value of the RNF sheet's last row (RNF-06) from "Media" to
# "ALTA".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Fix RNF row 7 (RNF-06) priority: "Media" -> "ALTA"
# ---------------------------------------------------------------------------
$rnf = $wb.Worksheets.Item("RNF")
$rnf.Range("E7").Value = "ALTA"

# ---------------------------------------------------------------------------
# 2. New sheet "Metricas" - defect density per module
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$metricas = $wb.Worksheets.Add($null, $lastSheet)
$metricas.Name = "Metricas"

$metricas.Range("A1").Value = "Modulo"
$metricas.Range("B1").Value = "Densidad de defectos"
$metricas.Range("C1").Value = "Numero de defectos"
$metricas.Range("D1").Value = "tamaño del modulo"

$metricas.Range("A2").Value = "Reportes"
$metricas.Range("C2").Value = 20
$metricas.Range("D2").Value = 1000
$metricas.Range("B2").Formula = "=`$C2/`$D2"

$metricas.Range("A3").Value = "Login"
$metricas.Range("C3").Value = 2
$metricas.Range("D3").Value = 500
$metricas.Range("B3").Formula = "=`$C3/`$D3"

# Highlight the module-name column (header + values) in yellow
$metricas.Range("A1:B1").Interior.Color = 65535
$metricas.Range("A2").Interior.Color = 65535
$metricas.Range("A3").Interior.Color = 65535

$metricas.Range("A1:D3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. New sheet "MTBF" - Mean Time Between Failures
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$mtbf = $wb.Worksheets.Add($null, $lastSheet)
$mtbf.Name = "MTBF"

$mtbf.Range("A1").Value = "MTBF"
$mtbf.Range("B1").Value = "Tiempo total Operación"
$mtbf.Range("C1").Value = "Numero de falla"

$mtbf.Range("B2").Value = 12
$mtbf.Range("C2").Value = 3
$mtbf.Range("A2").Formula = "=`$B2/`$C2"

$mtbf.Range("A1").Interior.Color = 65535

$mtbf.Range("A1:C2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. New sheet "MTTR" - Mean Time To Repair
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$mttr = $wb.Worksheets.Add($null, $lastSheet)
$mttr.Name = "MTTR"

$mttr.Range("A1").Value = "MTTR"
$mttr.Range("B1").Value = "Tiempo Total de reparación"
$mttr.Range("C1").Value = "Numero de fallas"

$mttr.Range("B2").Value = 40
$mttr.Range("C2").Value = 1
$mttr.Range("B3").Value = 20
$mttr.Range("C3").Value = 2
$mttr.Range("B4").Value = 60
$mttr.Range("C4").Value = 3
$mttr.Range("B5").Value = 30
$mttr.Range("C5").Value = 4
$mttr.Range("A2").Formula = "=(SUM(`$B2:`$B5)/`$C5)"

# Touch A3:A5 so the (empty) cells are materialised, matching the filled-down
# formatting of column A in the authored workbook
$mttr.Range("A3:A5").Font.Size = 11

# ---------------------------------------------------------------------------
# 5. Restore focus to the RNF sheet (matches the workbook's saved view state)
# ---------------------------------------------------------------------------
$rnf.Activate() | Out-Null
$rnf.Range("E7").Select() | Out-Null
